# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Only column G ("K") values change for the data rows (rows 2-44).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 1
    4  = 3
    5  = 1
    6  = 2
    7  = 0
    8  = 3
    9  = 2
    10 = 1
    11 = 2
    12 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 2
    17 = 0
    18 = 1
    19 = 2
    20 = 0
    21 = 1
    22 = 2
    23 = 0
    24 = 2
    25 = 2
    26 = 1
    27 = 2
    28 = 1
    29 = 2
    30 = 0
    31 = 1
    32 = 2
    33 = 1
    34 = 1
    35 = 1
    36 = 1
    37 = 2
    38 = 2
    39 = 2
    40 = 0
    41 = 2
    42 = 0
    43 = 1
    44 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
